$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "331.66"
Set-TextValue $ws "E2" "0.78%"
Set-TextValue $ws "D3" "45.45"
Set-TextValue $ws "E3" "3.42%"
Set-TextValue $ws "D4" "5.615"
Set-TextValue $ws "E4" "2.07%"
Set-TextValue $ws "E5" "4.54%"
Set-TextValue $ws "D6" "2.074"
Set-TextValue $ws "E6" "4.44%"
Set-TextValue $ws "D7" "0.9621"
Set-TextValue $ws "E7" "1.30%"
Set-TextValue $ws "D9" "0.1171"
Set-TextValue $ws "E9" "5.05%"
Set-TextValue $ws "D10" "0.1924"
Set-TextValue $ws "E10" "1.96%"
Set-TextValue $ws "D11" "10.41"
Set-TextValue $ws "E11" "-2.20%"
Set-TextValue $ws "D12" "0.09855"
Set-TextValue $ws "E12" "-1.84%"
Set-TextValue $ws "D13" "0.04618"
Set-TextValue $ws "E13" "-2.63%"
Set-TextValue $ws "D14" "0.1061"
Set-TextValue $ws "E14" "-0.22%"
Set-TextValue $ws "D15" "0.001293"
Set-TextValue $ws "E15" "2.23%"
Set-TextValue $ws "D16" "0.006118"
Set-TextValue $ws "E16" "2.53%"
Set-TextValue $ws "D17" "3.378"
Set-TextValue $ws "E17" "0.33%"
Set-TextValue $ws "D18" "4.442"
Set-TextValue $ws "E18" "1.49%"
Set-TextValue $ws "D19" "0.3341"
Set-TextValue $ws "E19" "-3.72%"
Set-TextValue $ws "D20" "0.1394"
Set-TextValue $ws "E20" "-1.89%"
Set-TextValue $ws "D21" "0.2656"
Set-TextValue $ws "E21" "2.54%"
Set-TextValue $ws "D22" "0.04188"
Set-TextValue $ws "E22" "2.37%"
Set-TextValue $ws "D24" "0.004567"
Set-TextValue $ws "E24" "7.14%"
Set-TextValue $ws "E25" "8.47%"
Set-TextValue $ws "D26" "0.0003752"
Set-TextValue $ws "D38" "0.02706"
Set-TextValue $ws "E38" "5.62%"
Set-TextValue $ws "D39" "0.05759"
Set-TextValue $ws "E39" "1.74%"
Set-TextValue $ws "D40" "0.007824"
Set-TextValue $ws "E40" "3.63%"
Set-TextValue $ws "D41" "0.1433"
Set-TextValue $ws "E41" "2.58%"
Set-TextValue $ws "D42" "0.007305"
Set-TextValue $ws "E42" "-1.37%"
Set-TextValue $ws "D43" "0.002016"
Set-TextValue $ws "E43" "-0.12%"
Set-TextValue $ws "D44" "0.009140"
Set-TextValue $ws "E44" "9.56%"
Set-TextValue $ws "D45" "0.3543"
Set-TextValue $ws "D46" "0.00007123"
Set-TextValue $ws "E46" "-0.26%"
Set-TextValue $ws "E47" "0.13%"
Set-TextValue $ws "D48" "0.0005819"
Set-TextValue $ws "E48" "0.13%"
Set-TextValue $ws "D49" "0.003493"
Set-TextValue $ws "E49" "-7.50%"
Set-TextValue $ws "D50" "0.003509"
Set-TextValue $ws "E50" "-0.73%"
Set-TextValue $ws "E51" "0.13%"
